# Generate Report for Handback
# The f6c6938f-6784-4149-ac3a-94204185f2d8.md file has now been handed
# back (in sync with en-US). Update its Status on every sheet, stamp the
# new "Latest Handback DateTime" on the zh-cn / de-de detail sheets, and
# clear the old "Error Detail" message now that it is no longer stale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for the f6c6938f row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# --- zh-cn detail sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("K3").Value = "2016-08-23 16:52:51"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.83

# --- de-de detail sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("K3").Value = "2016-08-23 16:52:58"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.83
